# Amélioration de l'inscription à la newsletter
# Adds a 6th test sheet ("Test CU 6 S'inscrire à la news") and tweaks
# sheet 1 / sheet 5 selections + styling to match.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Update "Test CU 5 Contacter l'assoc": selection moves + B7 turns
#    into the red "Ko" style (it was plain/unstyled before).
# ---------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("Test CU 5 Contacter l'assoc")
$ws5.Activate()
$ws5.Range("B7").Font.Color = 255        # BGR 255 == pure red == FFFF0000 (the "Ko" font)
$ws5.Range("B7").Select()

# ---------------------------------------------------------------------
# 2. Add the new sheet at the end of the workbook, same layout/style
#    pattern as the other "Test CU" sheets.
# ---------------------------------------------------------------------
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$ws6 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws6.Name = "Test CU 6 S'inscrire à la news"

# Header row
$ws6.Range("A1").Value = "Nom du test"
$ws6.Range("B1").Value = "Résultat du test"
$ws6.Range("C1").Value = "Observations"
$ws6.Range("A1:C1").Font.Bold = $true

# Body rows (A = step description, B = Ok result). Shared-string table
# order follows first-use order, and row 7's sentence was authored before
# rows 2-4's, so write it first to land on the right <si> index.
$ws6.Range("A7").Value = "6. Le membre reçoit un courriel qu'il est bien inscrit à la lettre d'information"
$ws6.Range("A2").Value = "1. Le membre se rend sur la page d'inscription à la newsletter"
$ws6.Range("A3").Value = "2. Le site lui affiche l'interface d'inscription à la newsletter"
$ws6.Range("A4").Value = "3. Le membre rentre son email et valide"
$ws6.Range("A5").Value = "4. Le site vérifie les informations saisies par l'internaute"
$ws6.Range("A6").Value = "5. Le site informe l'internaute de la véracité des informations saisies"

$ws6.Range("A2:A7").Font.Color = 12611584   # BGR == FF0070C0 blue, same as the other sheets
$ws6.Range("A2:A7").VerticalAlignment = -4108 # xlCenter

$ws6.Range("B2").Value = "Ok"
$ws6.Range("B3").Value = "Ok"
$ws6.Range("B4").Value = "Ok"
$ws6.Range("B5").Value = "Ok"
$ws6.Range("B6").Value = "Ok"
$ws6.Range("B7").Value = "Ok"
$ws6.Range("B2:B7").Font.Color = 5287936    # BGR == FF00B050 green, the "Ok" font

# Column widths, matching the bestFit pattern used on the other sheets
$ws6.Columns.Item(1).ColumnWidth = 73.29
$ws6.Columns.Item(2).ColumnWidth = 13.79
$ws6.Columns.Item(3).ColumnWidth = 11.88

# Page setup, same as sheets 1-4
$ws6.PageSetup.PaperSize = 9
$ws6.PageSetup.Orientation = 1

# Selection + active cell for the new sheet
$ws6.Range("A10").Select()

# ---------------------------------------------------------------------
# 3. Leave the new sheet active/selected (tabSelected + activeTab),
#    which also drops tabSelected from sheet 1.
# ---------------------------------------------------------------------
$ws6.Activate()
$win = $excel.ActiveWindow
$win.FirstSheet = 2
